$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.420.52"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "2.239.24"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +1.03%  "
$ws.Range("D5").Value = "'307.60"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'94.26"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'34.58"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").Value = "'0.0802"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "'7.19"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.228.40"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.832"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "'13.55"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "44.111.90"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "0.0₃0953"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "'6.34"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "'11.87"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").Value = "'65.46"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").Value = "'237.52"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").Value = "'1.97"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'2.22"
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("D27").Value = "'37.92"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").Value = "'9.76"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'19.90"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'153.13"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").Value = "'3.09"
$ws.Range("E34").Value = "  -3.37%  "
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("D37").Value = "'1.80"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("D38").Value = "'14.91"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").Value = "'3.38"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "'3.75"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "1.787.94"
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("D45").Value = "'78.47"
$ws.Range("E45").Value = "  -7.81%  "
$ws.Range("D46").Value = "'70.06"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'4.89"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'98.27"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("E49").Value = "  +4.87%  "
$ws.Range("D50").Value = "'8.06"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").Value = "'54.35"
$ws.Range("E51").Value = "  +0.18%  "
